# Update with latest cht-conf changes again and remove NO_LABEL
#
# The "survey" sheet's C3 cell ("NO_LABEL") is no longer needed for the
# begin_group/inputs row, so its contents are cleared. This also causes
# the now-unused "NO_LABEL" shared string to be dropped from the workbook
# on save (sharedStrings count/uniqueCount both shrink by one), which in
# turn shifts every other <v> shared-string index down by one - all
# handled automatically by the engine, the cell text values themselves
# are unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Remove the "NO_LABEL" text that used to live in C3.
$ws.Range("C3").ClearContents()

# Re-point the bottom-right pane's active cell/selection at C3 (previously
# C17), matching where the author's cursor ended up after editing.
$ws.Range("C3").Select()

Write-Output "Cleared survey!C3 (NO_LABEL) and updated selection to C3"
